$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add drilling dialog worked example under user 4:
#   E41 = drill point angle (deg), E42 = diameter -> E45 computes the
#   additional "tip" depth needed for full-depth / pecking drilling cycles.
$ws.Range("E41").Value = 59
$ws.Range("E42").Value = 2.1
$ws.Range("E45").Formula = "=E42/TAN(RADIANS(E41))"
$ws.Range("E47").Value = 1.2618073000000001

# Scroll the view down to the new block and leave the result cell selected,
# matching the author's final cursor position.
$ws.Range("A19").Select()
$excel.ActiveWindow.ScrollRow = 19
$ws.Range("E45").Select()
